# "maj - 22h30 a 01-47" : two mission-tracker updates were pasted in,
# one around 22h30 (still "today" = row 36) and one around 01h47 after
# midnight (new "today" = row 35, which stops being the hidden placeholder
# row and becomes the live day).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# --- cosmetic: move the active cell/selection ---
$ws.Range("X36").Select()

# --- Row 35 becomes the live "today" row: unhide it ---
$ws.Rows.Item(35).Hidden = $false

# --- Row 35: new snapshot values (01h47 update) ---
# C35 is left alone on purpose - it keeps its "=C36" formula and will
# recompute once C36 is updated below.
$ws.Range("D35").Value = 37201
$ws.Range("G35").Value = 33
$ws.Range("H35").Value = 39842
$ws.Range("K35").Value = 48
$ws.Range("L35").Value = 53751
$ws.Range("O35").Value = 69
$ws.Range("P35").Value = 47141
$ws.Range("S35").Value = 19
$ws.Range("T35").Value = 11283
$ws.Range("W35").Value = 49
$ws.Range("X35").Value = 46200
$ws.Range("AA35").Value = 46
$ws.Range("AB35").Value = 61270

# --- Row 36: updated snapshot values (22h30 update, still same day) ---
$ws.Range("C36").Value = 56
$ws.Range("D36").Value = 37200
$ws.Range("G36").Value = 30
$ws.Range("H36").Value = 37687
$ws.Range("K36").Value = 47
$ws.Range("L36").Value = 51770

# --- Row 36 now gets the rolling projection columns (AG:AJ) that the
# older, already-finalized rows (37, 38, ...) already have ---
$ws.Range("AG36").Value = 16592
$ws.Range("AH36").Formula = "=AG36-AG37"
$ws.Range("AI36").Formula = "=AE36/`$AH36"
$ws.Range("AJ36").Formula = "=AF36/`$AH36"
